# Applies the "split checks_and_parents / yield_checks" edit to check_genotypes.xlsx
$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "checks_and_parents"
$ws1.Range("A1:A15").Select() | Out-Null

# --- Add a new sheet for just the yield checks (no parent lines) ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "yield_checks"

$values = @("genotype","Dilday","Dunphy","N.C. Raleigh","Osage","Roy","NC-Dunphy","NC-Dilday","NC-Raleigh","Ellis","N10-687")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = $values[$i]
}

# Make the new sheet the active / selected one, matching tabSelected in the diff
$ws2.Select() | Out-Null
$ws2.Range("C10").Select() | Out-Null

$wb.Save()
